# Auto-generated edit script: applies cached-value updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# per the scheduled-runner price/profit recompute.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4663.3335
$ws.Range("J13").Value = 1995
$ws.Range("L13").Value = 1995
$ws.Range("N13").Value = -2333

$ws.Range("H17").Value = 1853.625
$ws.Range("J17").Value = 2227.158
$ws.Range("L17").Value = 6681.474
$ws.Range("N17").Value = -7017.474

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = $null

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = $null

$ws.Range("H96").Value = 388.6
$ws.Range("I96").Value = 279.5
$ws.Range("K96").Value = 838.5
$ws.Range("M96").Value = 534.5

$ws.Range("H100").Value = 3478.75
$ws.Range("I100").Value = 1218.125
$ws.Range("K100").Value = 1218.125
$ws.Range("M100").Value = -677.125

$ws.Range("H104").Value = 1131.6666
$ws.Range("I104").Value = 1131.6666
$ws.Range("K104").Value = 3394.9998
$ws.Range("M104").Value = -1647.9998

$ws.Range("H132").Value = 17919.066
$ws.Range("I132").Value = 21348.6
$ws.Range("K132").Value = 64045.8
$ws.Range("M132").Value = -61515.8

$ws.Range("H138").Value = 3382.5833
$ws.Range("I138").Value = 1616.3334
$ws.Range("K138").Value = 4849.0002
$ws.Range("M138").Value = 290.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1600
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = $null

$ws.Range("H32").Value = 4399.077
$ws.Range("I32").Value = 4399.077
$ws.Range("K32").Value = 4399.077
$ws.Range("M32").Value = -4112.077

$ws.Range("H61").Value = 6417.4443
$ws.Range("I61").Value = 6233.1665
$ws.Range("K61").Value = 6233.1665
$ws.Range("M61").Value = -6021.1665

$ws.Range("H88").Value = 1917.5714
$ws.Range("I88").Value = 1400
$ws.Range("J88").Value = 2124.6
$ws.Range("K88").Value = 1400
$ws.Range("L88").Value = 2124.6
$ws.Range("M88").Value = -994
$ws.Range("N88").Value = -2936.6

$ws.Range("H91").Value = 1917.5714
$ws.Range("I91").Value = 1400
$ws.Range("J91").Value = 2124.6
$ws.Range("K91").Value = 1400
$ws.Range("L91").Value = 2124.6
$ws.Range("M91").Value = 4
$ws.Range("N91").Value = -4932.6

$ws.Range("H132").Value = 2660.6875
$ws.Range("I132").Value = 1338.3636
$ws.Range("K132").Value = 4015.0908
$ws.Range("M132").Value = -1485.0908

$ws.Range("H136").Value = 6417.4443
$ws.Range("I136").Value = 6233.1665
$ws.Range("K136").Value = 18699.4995
$ws.Range("M136").Value = -16149.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4420.1816
$ws.Range("I86").Value = 915.5
$ws.Range("K86").Value = 915.5
$ws.Range("M86").Value = 207.5

$ws.Range("H89").Value = 4420.1816
$ws.Range("I89").Value = 915.5
$ws.Range("K89").Value = 4577.5
$ws.Range("M89").Value = 1038.5

$ws.Range("H134").Value = 3769.3333
$ws.Range("I134").Value = 2988.3845
$ws.Range("J134").Value = 5799.8
$ws.Range("K134").Value = 8965.1535
$ws.Range("L134").Value = 17399.4
$ws.Range("M134").Value = -6430.1535
$ws.Range("N134").Value = -22469.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2000155.4
$ws.Range("I3").Value = 3333467.2
$ws.Range("J3").Value = 187.5
$ws.Range("K3").Value = 3333467.2
$ws.Range("L3").Value = 187.5
$ws.Range("M3").Value = -3333354.2
$ws.Range("N3").Value = -413.5

$ws.Range("H13").Value = 400
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -261
$ws.Range("N13").Value = $null

$ws.Range("H31").Value = 5709.2144
$ws.Range("I31").Value = 1208.2
$ws.Range("J31").Value = 8209.777
$ws.Range("K31").Value = 1208.2
$ws.Range("L31").Value = 8209.777
$ws.Range("M31").Value = -913.2
$ws.Range("N31").Value = -8799.777

$ws.Range("H34").Value = 5709.2144
$ws.Range("I34").Value = 1208.2
$ws.Range("J34").Value = 8209.777
$ws.Range("K34").Value = 1208.2
$ws.Range("L34").Value = 8209.777
$ws.Range("M34").Value = -1006.2
$ws.Range("N34").Value = -8613.777

$ws.Range("H58").Value = 4014.8333
$ws.Range("I58").Value = 1524.75
$ws.Range("J58").Value = 8995
$ws.Range("K58").Value = 1524.75
$ws.Range("L58").Value = 8995
$ws.Range("M58").Value = -1321.75
$ws.Range("N58").Value = -9401

$ws.Range("H105").Value = 1639.625
$ws.Range("I105").Value = 1639.625
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1639.625
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 107.375
$ws.Range("N105").Value = $null

$ws.Range("H132").Value = 4171.533
$ws.Range("I132").Value = 3548.6667
$ws.Range("J132").Value = 6663
$ws.Range("K132").Value = 10646.0001
$ws.Range("L132").Value = 19989
$ws.Range("M132").Value = -8116.000100000001
$ws.Range("N132").Value = -25049

$ws.Range("H136").Value = 4014.8333
$ws.Range("I136").Value = 1524.75
$ws.Range("J136").Value = 8995
$ws.Range("K136").Value = 4574.25
$ws.Range("L136").Value = 26985
$ws.Range("M136").Value = -2024.25
$ws.Range("N136").Value = -32085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 222
$ws.Range("I86").Value = 205
$ws.Range("K86").Value = 615
$ws.Range("M86").Value = 571

$ws.Range("H89").Value = 222
$ws.Range("I89").Value = 205
$ws.Range("K89").Value = 1845
$ws.Range("M89").Value = 4083

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4854.3335
$ws.Range("I55").Value = 5700
$ws.Range("J55").Value = 4177.8
$ws.Range("K55").Value = 5700
$ws.Range("L55").Value = 4177.8
$ws.Range("M55").Value = -5373
$ws.Range("N55").Value = -4831.8

$ws.Range("H80").Value = 1789.8
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 2383
$ws.Range("K80").Value = 900
$ws.Range("L80").Value = 2383
$ws.Range("M80").Value = 98
$ws.Range("N80").Value = -4379

$ws.Range("H83").Value = 1789.8
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 2383
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 11915
$ws.Range("M83").Value = 492
$ws.Range("N83").Value = -21899

$ws.Range("H132").Value = 119117.11
$ws.Range("I132").Value = 119117.11
$ws.Range("K132").Value = 357351.33
$ws.Range("M132").Value = -354821.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6602.5386
$ws.Range("I68").Value = 4910
$ws.Range("J68").Value = 7660.375
$ws.Range("K68").Value = 4910
$ws.Range("L68").Value = 7660.375
$ws.Range("M68").Value = -4161
$ws.Range("N68").Value = -9158.375

$ws.Range("H71").Value = 6602.5386
$ws.Range("I71").Value = 4910
$ws.Range("J71").Value = 7660.375
$ws.Range("K71").Value = 24550
$ws.Range("L71").Value = 38301.875
$ws.Range("M71").Value = -20806
$ws.Range("N71").Value = -45789.875

$ws.Range("H82").Value = 2734
$ws.Range("I82").Value = 1149.5
$ws.Range("K82").Value = 1149.5
$ws.Range("M82").Value = -788.5

$ws.Range("H85").Value = 2734
$ws.Range("I85").Value = 1149.5
$ws.Range("K85").Value = 1149.5
$ws.Range("M85").Value = 98.5

$ws.Range("H136").Value = 4013.3333
$ws.Range("I136").Value = 3030.25
$ws.Range("K136").Value = 9090.75
$ws.Range("M136").Value = -6540.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4976
$ws.Range("I3").Value = 4666.6665
$ws.Range("K3").Value = 4666.6665
$ws.Range("M3").Value = -4552.6665

$ws.Range("H14").Value = 12501000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336

$ws.Range("H17").Value = 2418
$ws.Range("I17").Value = 2418
$ws.Range("K17").Value = 2418
$ws.Range("M17").Value = -2246

$ws.Range("H81").Value = 600
$ws.Range("I81").Value = 600
$ws.Range("K81").Value = 1200
$ws.Range("M81").Value = -139

$ws.Range("H84").Value = 600
$ws.Range("I84").Value = 600
$ws.Range("K84").Value = 6000
$ws.Range("M84").Value = -696

$ws.Range("H136").Value = 5495.0835
$ws.Range("I136").Value = 4280
$ws.Range("K136").Value = 12840
$ws.Range("M136").Value = -10290
